# Kilimanjaro Weekly Scoreboard - append the latest week's workout rows
# (rows 57-63) to the bottom of the Sheet1 data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column order: Participant, Date, Workout Type, Total Duration, Total Distance,
#               Total Elevation, Zone 1, Zone 2, Zone 3, Zone 4, Zone 5,
#               Workout Level, Week
$newRows = @(
    @("Eric",     45460, "Run",     12, 1.23, 43,  0,  0,  1, 7,  0, "Brave Leopard",  2),
    @("Eric",     45460, "Workout", 75, 0,    0,   2,  36, 28, 9, 0, "Brave Leopard",  2),
    @("Steven",   45460, "Walk",    31, 1.5,  36,  31, 0,  0,  0, 0, "Agile Antelope", 2),
    @("Jeremiah", 45460, "Run",     37, 3.74, 144, 1,  23, 10, 0, 0, "Agile Antelope", 2),
    @("Jeremiah", 45461, "Run",     48, 4.69, 69,  0,  29, 16, 0, 0, "Agile Antelope", 2),
    @("Matt",     45461, "Run",     37, 4,    272, 4,  28, 3,  0, 0, "Agile Antelope", 2),
    @("Matt",     45461, "Walk",    8,  0.29, 30,  8,  0,  0,  0, 0, "Agile Antelope", 2)
)

$startRow = 57
$endRow = $startRow + $newRows.Count - 1

# Pick up the existing date number format (column B already uses it, e.g. B2)
# so the new Date cells render/serialize with the same style instead of a
# brand-new number format being created.
$ws.Range("B2").Copy()
$ws.Range("B" + $startRow + ":B" + $endRow).PasteSpecial(-4122)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Match the author's final selection/active cell (bottom-right of new data).
$ws.Range("M" + $endRow).Select()
